$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Rewrite the "My idea was to create a Pokémon..." paragraph.
#    The original paragraph had many runs (one per repeated mention
#    of "Pokémon"). The edit trims it down to a much shorter
#    paragraph, keeping the leading "My idea was to create a " /
#    "Pokémon" runs intact and replacing everything from
#    "-centered program." onward.
# -----------------------------------------------------------------

$ideaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "My idea was to create a*") {
        $ideaPara = $candidate
        break
    }
}

if ($ideaPara -ne $null) {
    $pPr = $ideaPara.Range.ParagraphFormat
    $rng = $ideaPara.Range

    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5991D152" w14:textId="1B2BF576" w:rsidR="004C5AB1" w:rsidRDefault="004C5AB1"><w:r><w:t xml:space="preserve">My idea was to create a </w:t></w:r><w:r w:rsidRPr="004C5AB1"><w:t>Pok' + [char]0x00E9 + 'mon</w:t></w:r><w:r><w:t xml:space="preserve">-centered program. </w:t></w:r><w:r><w:t>The two maintainable tables I am going to use are Pok' + [char]0x00E9 + 'mon and Trainer. Tables for Type, Nature, and Region will be populated as they will be used in a dropdown for the Trainer and Pok' + [char]0x00E9 + 'mon maintenance forms.</w:t></w:r></w:p>'

    $rng.InsertXML($newParaXml)
}

# -----------------------------------------------------------------
# 2) Update the two business-rule bullet points.
# -----------------------------------------------------------------

$d.Content.Find.Execute("A battle can have two teams", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "A trainer can have up to 6 Pok" + [char]0x00E9 + "mon", 2)

$d.Content.Find.Execute("Both teams cannot have the same trainer", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "A trainer can only have one of each Pok" + [char]0x00E9 + "mon", 2)
